# Release 02 of "The Adventures of an Adventurer" + Release and Planing Docs
#
# This script updates the "Guppen-Planung" task table:
#   - marks several existing tasks as fully completed (Remain = 0)
#   - adds three brand-new tasks (rows 53-55) that were not tracked before
#   - moves the view/selection further down the sheet to where the new
#     work items live
#
# Columns: A=Task  B=Orig.Est.(h)  C=Curr.Est.(h)  D=Effort(h)
#          E=Remain(h) (=C-D)      F=Responsible   G=Completion(%) (=1-E/C)
#          H=Completed (=IF(G=100%,"checkmark","X"))  I=User Story ID
#          J=Datum   K=Notes

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15: "Bugfix Slime AI for Edges" -> fully done (Effort now equals Curr. Est.) ---
$ws.Cells.Item(15, 4).Value = 2

# --- Row 35: "Basic Level Selection..." -> fully done ---
$ws.Cells.Item(35, 4).Value = 5

# --- Row 50: "Add additional attack for first boss" -> fully done ---
$ws.Cells.Item(50, 4).Value = 10

# --- Row 51: "Implementing XML Reader for all Dialoges" -> scope grew, now fully done ---
$ws.Cells.Item(51, 3).Value = 4
$ws.Cells.Item(51, 4).Value = 4

# --- Row 52: "Landscape sprites/tiles" -> fully done ---
$ws.Cells.Item(52, 4).Value = 6

# --- Row 53 (new task): Porky(enemy) implementation in game ---
$ws.Cells.Item(53, 1).Value = "Porky(enemy) implementation in game"
$ws.Cells.Item(53, 2).Value = 5
$ws.Cells.Item(53, 3).Value = 5
$ws.Cells.Item(53, 4).Value = 2
$ws.Cells.Item(53, 6).Value = "Sascha"
$ws.Cells.Item(53, 9).Value = 13

# --- Row 54 (new task): Bat enemy implementation ---
$ws.Cells.Item(54, 1).Value = "Bat enemy implementation"
$ws.Cells.Item(54, 2).Value = 4
$ws.Cells.Item(54, 3).Value = 4
$ws.Cells.Item(54, 4).Value = 0
$ws.Cells.Item(54, 6).Value = "Cedric"
$ws.Cells.Item(54, 9).Value = 13

# --- Row 55 (new task): Implement first level (Basic Colliders, "Decoration") ---
$ws.Cells.Item(55, 1).Value = "Implement first level (Basic Colliders, ""Decoration"")"
$ws.Cells.Item(55, 2).Value = 3
$ws.Cells.Item(55, 3).Value = 3
$ws.Cells.Item(55, 6).Value = "Cedric, Sascha"
$ws.Cells.Item(55, 9).Value = 13

# --- Update the visible window/selection to scroll down toward the new rows ---
# (the sheet already has a frozen header pane; just move the active selection)
[void]$ws.Range("A55").Select()
